# Generate Report for Handback
# Refresh the timestamp values recorded on the handback-status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G, row 2)
$wsOverview.Range("G2").Value = "2016-08-31 13:23:24"

# zh-cn sheet: "Correspond Handoff Datetime" (column H) / "Correspond Handback DateTime" (column K)
$wsZhCn.Range("H2").Value = "2016-08-31 13:23:19"
$wsZhCn.Range("K2").Value = "2016-08-31 13:23:36"

# de-de sheet: "Correspond Handback DateTime" (column K)
$wsDeDe.Range("K2").Value = "2016-08-31 13:23:43"
